# Generate Report for Handoff
# Adds a new tracked file ("175c5e52-4b3b-4df7-b275-7406c7a6b704ooo....md") as
# row 3 on the Overview / zh-cn / de-de sheets, mirroring the pre-existing
# "8be37c0f-..." row but with a "Ready for handoff" status.

$wb = $excel.ActiveWorkbook

$newFileBare = "175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newFileE2e  = "e2e\175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$statusText  = "Ready for handoff"
$dedeDate    = "2016-09-03 08:31:33"
$zhcnDate    = "2016-09-03 08:31:29"
$zhcnXlf     = "175c5e52-4b3b-4df7-b275-7406c7a6b704oooooooooooooooooooooooooooooooooooooooo.0e2600ce9f4c4300e6d0a0806bca16914d3f1300.zh-cn.xlf"
$dedeXlf     = "175c5e52-4b3b-4df7-b275-7406c7a6b704oooooooooooooooooooooooooooooooooooooooo.0e2600ce9f4c4300e6d0a0806bca16914d3f1300.de-de.xlf"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e80a0ace2511f16f40f71139be338c78456157f1/e2e/" + $newFileBare

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFileBare
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $baseUrl, "", "", $newFileE2e) | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = $dedeDate
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $baseUrl, "", "", $newFileBare) | Out-Null
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = $zhcnXlf
$wsZhCn.Range("H3").Value = $zhcnDate
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = ""

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $baseUrl, "", "", $newFileBare) | Out-Null
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = $dedeXlf
$wsDeDe.Range("H3").Value = $dedeDate
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = ""

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
